$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/identifier-type"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"
$metadata.Range("B11").Value = "Extended set of identifier type codes from FHIR HL7 and extended LinuxForHealth Common Data Model resource identifiers"

$includeType = $wb.Worksheets.Item("Include from Identifier Type ")
$includeType.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/identifier-type"
